$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension-relevant data: rows 2..27, columns B..F
$data = @{
    2 = @('NSE:APOLLOTYRE', 'NSE:ASAL', $null, 'NSE:BHARATFORG', 'NSE:APOLLOTYRE')
    3 = @('NSE:BBETF0432', 'NSE:BHANDARI', $null, 'NSE:GODREJCP', 'NSE:M&M')
    4 = @('NSE:CESC', 'NSE:CSBBANK', $null, $null, 'NSE:MRF')
    5 = @('NSE:CHAMBLFERT', 'NSE:EMSLIMITED', $null, $null, $null)
    6 = @('NSE:CRAFTSMAN', 'NSE:GLOBAL', $null, $null, $null)
    7 = @('NSE:DBCORP', 'NSE:GODFRYPHLP', $null, $null, $null)
    8 = @('NSE:DEEPINDS', 'NSE:HARDWYN', $null, $null, $null)
    9 = @('NSE:GOACARBON', 'NSE:HLVLTD', $null, $null, $null)
    10 = @('NSE:HINDZINC', 'NSE:IRFC', $null, $null, $null)
    11 = @('NSE:JAGSNPHARM', 'NSE:JBMA', $null, $null, $null)
    12 = @('NSE:LUPIN', 'NSE:KRIDHANINF', $null, $null, $null)
    13 = @('NSE:MAHLIFE', 'NSE:LICNETFGSC', $null, $null, $null)
    14 = @('NSE:MUKTAARTS', 'NSE:MHLXMIRU', $null, $null, $null)
    15 = @('NSE:NSIL', 'NSE:MOTILALOFS', $null, $null, $null)
    16 = @('NSE:PCBL', 'NSE:NAGAFERT', $null, $null, $null)
    17 = @('NSE:SAMBHAAV', 'NSE:NELCAST', $null, $null, $null)
    18 = @($null, 'NSE:PANACEABIO', $null, $null, $null)
    19 = @($null, 'NSE:PPL', $null, $null, $null)
    20 = @($null, 'NSE:PRITI', $null, $null, $null)
    21 = @($null, 'NSE:RGL', $null, $null, $null)
    22 = @($null, 'NSE:RITES', $null, $null, $null)
    23 = @($null, 'NSE:ROTO', $null, $null, $null)
    24 = @($null, 'NSE:RRKABEL', $null, $null, $null)
    25 = @($null, 'NSE:RUSHIL', $null, $null, $null)
    26 = @($null, 'NSE:SAKSOFT', $null, $null, $null)
    27 = @($null, 'NSE:SALASAR', $null, $null, $null)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt 5; $i++) {
        $col = 2 + $i  # B=2 .. F=6
        $v = $vals[$i]
        if ($null -eq $v) {
            $ws.Cells.Item($r, $col).Value = ""
        } else {
            $ws.Cells.Item($r, $col).Value = $v
        }
    }
}

# Add index numbers and style for rows 23..27 (new rows)
# Column A uses the same style as existing index cells (e.g. A2), so copy
# the formatting from A2 first, then overwrite with the correct value.
$srcIndexCell = $ws.Cells.Item(2, 1)
for ($r = 23; $r -le 27; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $srcIndexCell.Copy($cell)
    $cell.Value = $r - 2
}
